$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamps for column Z (timestamp), rows 2-112, in order,
# corresponding to a re-run of the pcsmote logging process.
$timestamps = @(
    "2025-10-17T07:09:38.261747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.262747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.263747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.264747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.265747",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.344880",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.429645",
    "2025-10-17T07:09:38.440672",
    "2025-10-17T07:09:38.440672",
    "2025-10-17T07:09:38.440672",
    "2025-10-17T07:09:38.441185",
    "2025-10-17T07:09:38.529613",
    "2025-10-17T07:09:38.530614",
    "2025-10-17T07:09:38.530614",
    "2025-10-17T07:09:38.530614",
    "2025-10-17T07:09:38.530614",
    "2025-10-17T07:09:38.530614",
    "2025-10-17T07:09:38.531612",
    "2025-10-17T07:09:38.531612",
    "2025-10-17T07:09:38.531612",
    "2025-10-17T07:09:38.531612"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 26).Value = $timestamps[$i]
}

Write-Output "Updated $($timestamps.Length) timestamp cells in column Z"
